# Adding Test case for Authoring
# Target sheet is "Test Cases" (sheet1.xml / first tab in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# ---------------------------------------------------------------------
# Row 67 (existing "VerifyDraftPostTabDisplayForZeroDrafts" test case):
#   - D67 picks up the "left/top aligned, wrapped" style used throughout
#     the sheet (same look as D66/D68's new state).
#   - E67's result flips from PASS to SKIP.
# ---------------------------------------------------------------------
$ws.Range("D66").Copy()
$ws.Range("D67").PasteSpecial(-4122)
$ws.Range("D67").Value2 = "Y"
$ws.Range("E67").Value2 = "SKIP"

# ---------------------------------------------------------------------
# Row 68 (existing "DeleteDraftPostFromPostModal" test case):
#   - D68 gets the same style treatment as D67.
#   - E68's result flips from PASS to SKIP.
# ---------------------------------------------------------------------
$ws.Range("D66").Copy()
$ws.Range("D68").PasteSpecial(-4122)
$ws.Range("D68").Value2 = "Y"
$ws.Range("E68").Value2 = "SKIP"

# ---------------------------------------------------------------------
# New row 69: new Authoring test case
#   "VerifyPostTitleDisplayInDraftSection" / OPQA-1199
# Values are entered C, then B, then A (the order the strings were first
# authored in), so newly-created shared-string entries line up the same
# way they did when this test case was originally added.
# ---------------------------------------------------------------------
$ws.Range("C69").Value2 = "Verify that Draft Post section displays the post title`n[ Or Untitled if title is not provided] and  time draft was saved"
$ws.Range("B69").Value2 = "OPQA-1199"
$ws.Range("A69").Value2 = "VerifyPostTitleDisplayInDraftSection"
$ws.Range("D69").Value2 = "Y"
$ws.Range("E69").Value2 = "PASS"

# Match styling used by similar rows:
#  - A69/E69 look like A65/E65 (plain wrapped body cells)
#  - B69 looks like B64/B67 (TCID-style column with border)
#  - C69 looks like C64/C66 (wrapped description, taller row)
#  - D69 looks like D64 (Runmode column style)
$ws.Range("A65").Copy()
$ws.Range("A69").PasteSpecial(-4122)

$ws.Range("B64").Copy()
$ws.Range("B69").PasteSpecial(-4122)

$ws.Range("C64").Copy()
$ws.Range("C69").PasteSpecial(-4122)

$ws.Range("D64").Copy()
$ws.Range("D69").PasteSpecial(-4122)

$ws.Range("E65").Copy()
$ws.Range("E69").PasteSpecial(-4122)

# New row is a two-line entry like the other "description" rows (30pt tall).
$ws.Rows.Item(69).RowHeight = 30

# ---------------------------------------------------------------------
# View state: scroll position + active selection, as captured after
# authoring the new row.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D65").Select()
